# Applies the recorded changes to "Template - Projeto de Casos de teste.xlsx"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Plan1")
$ws2 = $wb.Worksheets.Item("Plan2")

# Order of writes mirrors the original edit session so new shared-string
# entries land in the same sequence as the recorded diff.

# New "Classes de Equivalência Exercitadas" column header on Plan2
$ws2.Range("F2").Value = "Classes de Equivalência Exercitadas"

# Plan1: tag the single-letter equivalence classes with numbers
$ws1.Range("C6").Value = "V (1)"
$ws1.Range("D6").Value = "I (2)"

# Plan2: example row values for the new columns
$ws2.Range("D3").Value = "9 2020"
$ws2.Range("E3").Value = "v"

# Plan1: remaining equivalence-class tags (C/D/E/F columns)
$ws1.Range("C10").Value = "I (3)"
$ws1.Range("D10").Value = "I (4)"
$ws1.Range("E10").Value = "I (5)"
$ws1.Range("F10").Value = "V (6)"

$ws1.Range("C14").Value = "I (7)"
$ws1.Range("D14").Value = "I (8)"
$ws1.Range("E14").Value = "I (9)"
$ws1.Range("F14").Value = "I (10)"

$ws1.Range("C15").Value = "I (11)"
$ws1.Range("D15").Value = "I (12)"
$ws1.Range("E15").Value = "I (13)"
$ws1.Range("F15").Value = "I (14)"

$ws1.Range("C16").Value = "I (15)"
$ws1.Range("D16").Value = "I (16)"
$ws1.Range("E16").Value = "I (17)"
$ws1.Range("F16").Value = "I (18)"

$ws1.Range("C17").Value = "I (19)"
$ws1.Range("D17").Value = "I (20)"
$ws1.Range("E17").Value = "I (21)"
$ws1.Range("F17").Value = "V (22)"

# Plan1: J column output tags
$ws1.Range("J6").Value = "365 dias (23)"
$ws1.Range("J7").Value = "366 dias (24)"
$ws1.Range("J8").Value = "Mensagem de ano inválido (25)"

$ws1.Range("J13").Value = "31 dias (26)"
$ws1.Range("J14").Value = "30 dias (27)"
$ws1.Range("J15").Value = "29 dias (28)"
$ws1.Range("J16").Value = "28 dias (29)"

$ws1.Range("J17").Value = "Mensagem de erro mês invalido (30)"
$ws1.Range("J18").Value = "Mensagem de erro ano inválido (31)"
$ws1.Range("J19").Value = "Mesagem de erro tudo inválido (32)"
$ws1.Range("J20").Value = "mensagem erro (33)"

# Plan2: last new value
$ws2.Range("C3").Value = 2
$ws2.Range("F3").Value = "1, 29"

# --- Column width tweaks ---
# Column J on Plan1 widened 26.5 -> 33 chars
$ws1.Columns.Item(10).ColumnWidth = 32.166666666666664

# Plan2 columns D,E widened slightly and F added (16 -> ~16.83, 15.16 -> ~16.33, new 29)
$ws2.Columns.Item(4).ColumnWidth = 16.0
$ws2.Columns.Item(5).ColumnWidth = 15.5
$ws2.Columns.Item(6).ColumnWidth = 28.166666666666664

# --- View / selection changes: active tab moves from Plan2 to Plan1 ---
$ws2.Activate() | Out-Null
$ws2.Range("F3").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F27").Select() | Out-Null
$excel.ActiveWindow.Zoom = 112
